$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to touch to remain plain text,
# matching the workbook author's original inline-string representation
# (several new values, e.g. "1.00", would otherwise be auto-converted to
# numbers by Excel's General number format and lose their trailing zeros).
$priceCells = @("D2", "D3", "D5", "D6", "D10", "D13", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D28", "D30", "D33", "D34", "D35", "D37", "D38", "D40", "D44", "D45", "D47", "D50")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Updated coin prices (column D)
$ws.Range("D2").Value = "76.132.76"
$ws.Range("D3").Value = "3.032.41"
$ws.Range("D5").Value = "197.01"
$ws.Range("D6").Value = "619.31"
$ws.Range("D10").Value = "3.030.84"
$ws.Range("D13").Value = "5.25"
$ws.Range("D14").Value = "3.587.98"
$ws.Range("D15").Value = "28.84"
$ws.Range("D16").Value = "76.095.36"
$ws.Range("D18").Value = "3.027.26"
$ws.Range("D19").Value = "13.42"
$ws.Range("D20").Value = "8.93"
$ws.Range("D21").Value = "381.05"
$ws.Range("D22").Value = "2.37"
$ws.Range("D24").Value = "72.55"
$ws.Range("D25").Value = "3.168.78"
$ws.Range("D28").Value = "9.75"
$ws.Range("D30").Value = "0.997"
$ws.Range("D33").Value = "492.87"
$ws.Range("D34").Value = "1.92"
$ws.Range("D35").Value = "1.00"
$ws.Range("D37").Value = "0.121"
$ws.Range("D38").Value = "162.02"
$ws.Range("D40").Value = "190.44"
$ws.Range("D44").Value = "5.12"
$ws.Range("D45").Value = "0.771"
$ws.Range("D47").Value = "41.45"
$ws.Range("D50").Value = "0.593"

# Updated 1h volume/change percentages (column E)
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("E3").Value = "  +3.98%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  +4.54%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +6.29%  "
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("E14").Value = "  +4.03%  "
$ws.Range("E15").Value = "  +3.67%  "
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("E18").Value = "  +3.50%  "
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("E20").Value = "  +2.68%  "
$ws.Range("E21").Value = "  +3.69%  "
$ws.Range("E22").Value = "  +5.97%  "
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +5.97%  "
$ws.Range("E36").Value = "  +2.68%  "
$ws.Range("E37").Value = "  +11.79%  "
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("E40").Value = "  +6.21%  "
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("E42").Value = "  -5.21%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  +4.68%  "
$ws.Range("E45").Value = "  +18.03%  "
$ws.Range("E46").Value = "  +5.91%  "
$ws.Range("E47").Value = "  +3.29%  "
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  +7.81%  "
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("E51").Value = "  +0.36%  "
